$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.751.55"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.312.62"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.70"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.48"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.13"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.92"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0783"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.72"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.671.45"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.322.89"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.698.43"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.11"
$ws.Range("E19").Value = "  -5.05%  "
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("E23").Value = "  +4.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.13"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("E28").Value = "  +14.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.94"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.12"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.45"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.48"
$ws.Range("E42").Value = "  +19.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.923.25"
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.04"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.08"
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.540.69"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.26"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.14"
$ws.Range("E51").Value = "  +1.72%  "
